$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Partha"
$ws.Range("B7").Value = "Saradhi"
$ws.Range("C7").Value = "Pune, Maharashtra, 412207"
$ws.Range("D7").Value = "886-766-2916"
$ws.Range("E7").Value = "paardhu@gmail.com"
$ws.Range("E7").Interior.Pattern = -4142
$ws.Range("F7").Value = 43
$ws.Range("G7").Value = 4
